$d = $word.ActiveDocument

# Locate the bullet paragraph "Utilize R programming..." which is the last
# bullet in the GRAIL / Senior Statistical Programmer Analyst section.
# Note: Paragraph.Index is NOT the 1-based ordinal within $d.Paragraphs, so
# we find the ordinal position by manually walking the collection instead.
$anchorIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Utilize R programming to analyze survival data and assess risk factors*") {
        $anchorIndex = $i
    }
}

$newBullets = @(
    "Developed and maintained Linux cluster environments, including the installation and management of bioinformatics tools for large-scale genomic analysis (Docker).",
    "Managed the installation and configuration of R packages across multiple environments, ensuring compatibility and stability.",
    "Stored, indexed, and transferred large-scale genomics and clinical trial data, utilizing high-performance storage solutions and cloud platforms (AWS)."
)

foreach ($text in $newBullets) {
    $cur = $d.Paragraphs.Item($anchorIndex)
    $cur.Range.InsertParagraphAfter()
    $anchorIndex = $anchorIndex + 1
    $newPara = $d.Paragraphs.Item($anchorIndex)
    $newPara.Range.Text = $text
}
